# Update cryptos list values per the scraped diff (prices / 1h volume %).
# Some "Price" values are plain decimal numbers as text (e.g. "214.50"); Excel's
# normal type-inference would silently store those as real numbers, which
# would not match the workbook's original inline-string (text) cells and
# would introduce float rounding artifacts. Forcing NumberFormat "@" (Text)
# before assigning those values keeps them as text, matching the source data.
# Values that already contain two dots (e.g. "25.926.05") or a percent sign
# are not parsed as numbers by Excel and need no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.926.05"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.633.96"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.51%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.50"
$ws.Range("E5").Value = "  +0.08%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.46%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0632"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"
$ws.Range("E10").Value = "  -0.04%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.28%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.859.67"
$ws.Range("E12").Value = "  +0.39%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.655.58"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.24"
$ws.Range("E14").Value = "  -0.28%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("E15").Value = "  -1.83%  "

# Row 16 - now ShibaInu (was Litecoin)
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₃0755"
$ws.Range("E16").Value = "  -0.82%  "

# Row 17 - now Litecoin (was ShibaInu)
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.74"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.910.96"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.42%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.88"
$ws.Range("E20").Value = "  +1.00%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.46%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +0.36%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  -0.35%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.39%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.19"
$ws.Range("E25").Value = "  +0.60%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.30%  "

# Row 27 - Stellar
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.125"
$ws.Range("E27").Value = "  +1.66%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.29%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.46"
$ws.Range("E29").Value = "  -0.11%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.00%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0498"
$ws.Range("E31").Value = "  +0.96%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  -0.79%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +0.02%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -0.08%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.92%  "

# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.138.02"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +1.32%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -1.15%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.60%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.60%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("E42").Value = "  +0.20%  "

# Row 43 - Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.20"
$ws.Range("E43").Value = "  -1.43%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -2.21%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.768.95"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +0.51%  "

# Row 47 - Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.27"
$ws.Range("E47").Value = "  +2.06%  "

# Row 48 - Cronos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0523"
$ws.Range("E48").Value = "  +2.40%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  +0.70%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.20%  "

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("E51").Value = "  +1.64%  "
